$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("empdata")

# Row 3: fill in the previously-empty "status" cell (E3) to match E2's message
$ws.Range("E3").Value = "Employee details saved successfully"

# Row 4: new employee record (fname, lname, email, mobile)
$ws.Range("A4").Value = "raj"
$ws.Range("B4").Value = "rajiv"
$ws.Range("C4").Value = "raj@hot.com"
# store mobile number as text (leading apostrophe keeps it a text value,
# matching the existing D3 cell which is also text)
$ws.Range("D4").Value = "'8564689656"

# Turn the new email address into a mailto hyperlink, like C2/C3
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:raj@hot.com")

# Hyperlinks.Add() re-styles the cell on its own; restore the same
# "hyperlink" cell format already used by C3 (border + hyperlink font)
$ws.Range("C3").Copy($ws.Range("C4"))
$ws.Range("C4").Value = "raj@hot.com"

